# PowerShell COM-interop script applying the scheduled-runner price refresh
# to Sheets/Behemoth_Profits.xlsx (workbook tabs ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Each row's currentAveragePrice* (H/I/J), Leve Price (K/L) and Profit (M/N)
# columns are refreshed in place with freshly-pulled market-board values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 4450.385
$ws.Range("I106").Value = 2643.6667
$ws.Range("K106").Value = 2643.6667
$ws.Range("M106").Value = -2012.6667
$ws.Range("H132").Value = 1503.2354
$ws.Range("I132").Value = 1315.9375
$ws.Range("K132").Value = 3947.8125
$ws.Range("M132").Value = -1417.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10643585
$ws.Range("I32").Value = 14287341
$ws.Range("K32").Value = 14287341
$ws.Range("M32").Value = -14287054
$ws.Range("H113").Value = 48198
$ws.Range("I113").Value = 36000
$ws.Range("J113").Value = 51247.5
$ws.Range("K113").Value = 36000
$ws.Range("L113").Value = 51247.5
$ws.Range("M113").Value = -31661
$ws.Range("N113").Value = -59925.5
$ws.Range("H132").Value = 4057.2942
$ws.Range("I132").Value = 4330.7334
$ws.Range("K132").Value = 12992.2002
$ws.Range("M132").Value = -10462.2002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3286.4
$ws.Range("I20").Value = 3816.5833
$ws.Range("J20").Value = 1165.6666
$ws.Range("K20").Value = 3816.5833
$ws.Range("L20").Value = 1165.6666
$ws.Range("M20").Value = -3569.5833
$ws.Range("N20").Value = -1659.6666
$ws.Range("H86").Value = 3012.2104
$ws.Range("I86").Value = 3181.0715
$ws.Range("J86").Value = 2539.4
$ws.Range("K86").Value = 3181.0715
$ws.Range("L86").Value = 2539.4
$ws.Range("M86").Value = -2058.0715
$ws.Range("N86").Value = -4785.4
$ws.Range("H89").Value = 3012.2104
$ws.Range("I89").Value = 3181.0715
$ws.Range("J89").Value = 2539.4
$ws.Range("K89").Value = 15905.3575
$ws.Range("L89").Value = 12697
$ws.Range("M89").Value = -10289.3575
$ws.Range("N89").Value = -23929

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 759.4
$ws.Range("J22").Value = 1002
$ws.Range("L22").Value = 1002
$ws.Range("N22").Value = -1702
$ws.Range("H31").Value = 780095.6
$ws.Range("J31").Value = 1168529.9
$ws.Range("L31").Value = 1168529.9
$ws.Range("N31").Value = -1169119.9
$ws.Range("H34").Value = 780095.6
$ws.Range("J34").Value = 1168529.9
$ws.Range("L34").Value = 1168529.9
$ws.Range("N34").Value = -1168933.9
$ws.Range("H41").Value = 44125.625
$ws.Range("I41").Value = 3000
$ws.Range("J41").Value = 50000.715
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 50000.715
$ws.Range("M41").Value = -2572
$ws.Range("N41").Value = -50856.715
$ws.Range("H51").Value = 52200
$ws.Range("I51").Value = 22500
$ws.Range("K51").Value = 22500
$ws.Range("M51").Value = -21764
$ws.Range("H61").Value = 52200
$ws.Range("I61").Value = 22500
$ws.Range("K61").Value = 22500
$ws.Range("M61").Value = -22152
$ws.Range("H68").Value = 62998
$ws.Range("J68").Value = 62998
$ws.Range("L68").Value = 62998
$ws.Range("N68").Value = -64496
$ws.Range("H71").Value = 62998
$ws.Range("J71").Value = 62998
$ws.Range("L71").Value = 188994
$ws.Range("N71").Value = -196482
$ws.Range("H86").Value = 6000.6665
$ws.Range("I86").Value = 5574.7
$ws.Range("J86").Value = 6533.125
$ws.Range("K86").Value = 5574.7
$ws.Range("L86").Value = 6533.125
$ws.Range("M86").Value = -4451.7
$ws.Range("N86").Value = -8779.125
$ws.Range("H89").Value = 6000.6665
$ws.Range("I89").Value = 5574.7
$ws.Range("J89").Value = 6533.125
$ws.Range("K89").Value = 27873.5
$ws.Range("L89").Value = 32665.625
$ws.Range("M89").Value = -22257.5
$ws.Range("N89").Value = -43897.625
$ws.Range("H105").Value = 1379.1786
$ws.Range("I105").Value = 1314.8182
$ws.Range("K105").Value = 1314.8182
$ws.Range("M105").Value = 432.1818000000001
$ws.Range("H132").Value = 2371.5454
$ws.Range("I132").Value = 2120.0557
$ws.Range("J132").Value = 3503.25
$ws.Range("K132").Value = 6360.1671
$ws.Range("L132").Value = 10509.75
$ws.Range("M132").Value = -3830.1671
$ws.Range("N132").Value = -15569.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 327.2
$ws.Range("I34").Value = 178.66667
$ws.Range("K34").Value = 536.00001
$ws.Range("M34").Value = -452.00001
$ws.Range("H131").Value = 5963.5312
$ws.Range("J131").Value = 5963.5312
$ws.Range("L131").Value = 17890.5936
$ws.Range("N131").Value = -27970.5936

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13450.546
$ws.Range("I70").Value = 17422.572
$ws.Range("J70").Value = 6499.5
$ws.Range("K70").Value = 17422.572
$ws.Range("L70").Value = 6499.5
$ws.Range("M70").Value = -17152.572
$ws.Range("N70").Value = -7039.5
$ws.Range("H73").Value = 13450.546
$ws.Range("I73").Value = 17422.572
$ws.Range("J73").Value = 6499.5
$ws.Range("K73").Value = 17422.572
$ws.Range("L73").Value = 6499.5
$ws.Range("M73").Value = -16486.572
$ws.Range("N73").Value = -8371.5
$ws.Range("H92").Value = 40750
$ws.Range("J92").Value = 40750
$ws.Range("L92").Value = 40750
$ws.Range("N92").Value = -44494
$ws.Range("H93").Value = 60000
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("M93").Value = -63744
$ws.Range("H95").Value = 53468.8
$ws.Range("J95").Value = 53468.8
$ws.Range("L95").Value = 53468.8
$ws.Range("N95").Value = -58960.8
$ws.Range("H108").Value = 99869.5
$ws.Range("J108").Value = 99869.5
$ws.Range("L108").Value = 99869.5
$ws.Range("N108").Value = -107549.5
$ws.Range("H109").Value = 45118.25
$ws.Range("J109").Value = 45118.25
$ws.Range("L109").Value = 45118.25
$ws.Range("N109").Value = -47198.25
$ws.Range("H110").Value = 96504.75
$ws.Range("J110").Value = 96504.75
$ws.Range("L110").Value = 96504.75
$ws.Range("N110").Value = -104684.75
$ws.Range("H128").Value = 116990
$ws.Range("J128").Value = 116990
$ws.Range("L128").Value = 116990
$ws.Range("N128").Value = -126950

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1581.8334
$ws.Range("I16").Value = 1581.8334
$ws.Range("K16").Value = 1581.8334
$ws.Range("M16").Value = -1411.8334
$ws.Range("H40").Value = 4899.2
$ws.Range("J40").Value = 5536.364
$ws.Range("L40").Value = 5536.364
$ws.Range("N40").Value = -5808.364
$ws.Range("H118").Value = 120500
$ws.Range("J118").Value = 120500
$ws.Range("L118").Value = 120500
$ws.Range("N118").Value = -123814
$ws.Range("H136").Value = 52930.03
$ws.Range("I136").Value = 7468.5264
$ws.Range("J136").Value = 119373.766
$ws.Range("K136").Value = 22405.5792
$ws.Range("L136").Value = 358121.298
$ws.Range("M136").Value = -19855.5792
$ws.Range("N136").Value = -363221.298
